$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit new (longer) text (bestFit-equivalent ~37.09 chars)
$ws.Columns.Item(2).ColumnWidth = 36.15

# Add row 51: rear driver
$ws.Range("A51").Formula = "=A50+1"
$ws.Range("B51").Value = "REV1 suspension assembly rear driver"

# Add row 52: rear passenger
$ws.Range("A52").Formula = "=A51+1"
$ws.Range("B52").Value = "REV1 suspension assembly rear passenger"

# Update selection to reflect new active cell
$ws.Range("B52").Select()
